# New crime data collected - weekly CompStat refresh (week of 4/7/2025 - 4/13/2025)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header text (Volume/Number + reporting week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  15"
$ws.Range("C9").Value = "Report Covering the Week  4/7/2025  Through  4/13/2025"

# --- Column E width bestFit adjustment ---
$ws.Columns.Item(5).ColumnWidth = 7.433768

# --- Stable style-13 (text) source cells used as format donors: D14 ("0") / E14 ("***.*") ---

# --- Row 14 ---
$ws.Cells.Item(14,3).Value = 1
$ws.Cells.Item(14,3).NumberFormat = '#,##0'
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,6).NumberFormat = '#,##0'
$ws.Cells.Item(14,9).Value = 2
$ws.Cells.Item(14,11).Value = 100
$ws.Cells.Item(14,12).Value = -33.333333333333
$ws.Cells.Item(14,13).Value = 100
$ws.Cells.Item(14,14).Value = -75

# --- Row 15 ---
$ws.Cells.Item(15,4).Value = "'0"
$ws.Cells.Item(14, 4).Copy() | Out-Null
$ws.Cells.Item(15,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,5).Value = "***.*"
$ws.Cells.Item(14, 5).Copy() | Out-Null
$ws.Cells.Item(15,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15,9).Value = 13
$ws.Cells.Item(15,11).Value = 225
$ws.Cells.Item(15,12).Value = 550
$ws.Cells.Item(15,14).Value = 116.666666666667

# --- Row 16 ---
$ws.Cells.Item(16,4).Value = 4
$ws.Cells.Item(16,5).Value = -25
$ws.Cells.Item(16,6).Value = 10
$ws.Cells.Item(16,8).Value = 0
$ws.Cells.Item(16,9).Value = 31
$ws.Cells.Item(16,10).Value = 44
$ws.Cells.Item(16,11).Value = -29.545454545454
$ws.Cells.Item(16,12).Value = -31.111111111111
$ws.Cells.Item(16,13).Value = -56.338028169014
$ws.Cells.Item(16,14).Value = -89.455782312925

# --- Row 17 ---
$ws.Cells.Item(17,3).Value = 4
$ws.Cells.Item(17,4).Value = 7
$ws.Cells.Item(17,5).Value = -42.857142857142
$ws.Cells.Item(17,6).Value = 29
$ws.Cells.Item(17,7).Value = 19
$ws.Cells.Item(17,8).Value = 52.631578947368
$ws.Cells.Item(17,9).Value = 90
$ws.Cells.Item(17,10).Value = 91
$ws.Cells.Item(17,11).Value = -1.098901098901
$ws.Cells.Item(17,12).Value = -6.25
$ws.Cells.Item(17,13).Value = 157.142857142857
$ws.Cells.Item(17,14).Value = -55

# --- Row 18 ---
$ws.Cells.Item(18,4).Value = 2
$ws.Cells.Item(18,4).NumberFormat = '#,##0'
$ws.Cells.Item(18,5).Value = -100
$ws.Cells.Item(18,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(18,6).Value = 3
$ws.Cells.Item(18,8).Value = -25
$ws.Cells.Item(18,9).Value = 28
$ws.Cells.Item(18,10).Value = 22
$ws.Cells.Item(18,11).Value = 27.272727272727
$ws.Cells.Item(18,12).Value = -3.448275862068
$ws.Cells.Item(18,13).Value = -24.324324324324
$ws.Cells.Item(18,14).Value = -90.410958904109

# --- Row 19 ---
$ws.Cells.Item(19,3).Value = 6
$ws.Cells.Item(19,4).Value = 6
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,7).Value = 23
$ws.Cells.Item(19,8).Value = -30.434782608695
$ws.Cells.Item(19,9).Value = 70
$ws.Cells.Item(19,10).Value = 86
$ws.Cells.Item(19,11).Value = -18.604651162790
$ws.Cells.Item(19,12).Value = -52.380952380952
$ws.Cells.Item(19,13).Value = -41.666666666666
$ws.Cells.Item(19,14).Value = -54.248366013071

# --- Row 20 ---
$ws.Cells.Item(20,3).Value = 7
$ws.Cells.Item(20,4).Value = 2
$ws.Cells.Item(20,5).Value = 250
$ws.Cells.Item(20,6).Value = 14
$ws.Cells.Item(20,7).Value = 5
$ws.Cells.Item(20,8).Value = 180
$ws.Cells.Item(20,9).Value = 28
$ws.Cells.Item(20,10).Value = 25
$ws.Cells.Item(20,11).Value = 12
$ws.Cells.Item(20,12).Value = 21.739130434782
$ws.Cells.Item(20,13).Value = -12.5
$ws.Cells.Item(20,14).Value = -90.604026845637

# --- Row 21 ---
$ws.Cells.Item(21,3).Value = 22
$ws.Cells.Item(21,4).Value = 21
$ws.Cells.Item(21,5).Value = 4.761904761904
$ws.Cells.Item(21,6).Value = 75
$ws.Cells.Item(21,7).Value = 63
$ws.Cells.Item(21,8).Value = 19.047619047619
$ws.Cells.Item(21,9).Value = 262
$ws.Cells.Item(21,10).Value = 273
$ws.Cells.Item(21,11).Value = -4.029304029304
$ws.Cells.Item(21,12).Value = -24.057971014492
$ws.Cells.Item(21,13).Value = -11.486486486486
$ws.Cells.Item(21,14).Value = -79.056754596322

# --- Row 22 ---
$ws.Cells.Item(22,3).Value = 2
$ws.Cells.Item(22,3).NumberFormat = '#,##0'
$ws.Cells.Item(22,6).Value = 4
$ws.Cells.Item(22,7).Value = "'0"
$ws.Cells.Item(14, 4).Copy() | Out-Null
$ws.Cells.Item(22,7).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,8).Value = "***.*"
$ws.Cells.Item(14, 5).Copy() | Out-Null
$ws.Cells.Item(22,8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(22,9).Value = 7
$ws.Cells.Item(22,11).Value = 16.666666666666
$ws.Cells.Item(22,12).Value = 133.333333333333
$ws.Cells.Item(22,13).Value = -46.153846153846

# --- Row 23 ---
$ws.Cells.Item(23,3).Value = 1
$ws.Cells.Item(23,4).Value = 1
$ws.Cells.Item(23,5).Value = 0
$ws.Cells.Item(23,6).Value = 8
$ws.Cells.Item(23,7).Value = 5
$ws.Cells.Item(23,8).Value = 60
$ws.Cells.Item(23,9).Value = 35
$ws.Cells.Item(23,10).Value = 43
$ws.Cells.Item(23,11).Value = -18.604651162790
$ws.Cells.Item(23,12).Value = 0
$ws.Cells.Item(23,13).Value = 66.666666666666

# --- Row 24 ---
$ws.Cells.Item(24,3).Value = 22
$ws.Cells.Item(24,4).Value = 12
$ws.Cells.Item(24,5).Value = 83.333333333333
$ws.Cells.Item(24,6).Value = 85
$ws.Cells.Item(24,7).Value = 73
$ws.Cells.Item(24,8).Value = 16.438356164383
$ws.Cells.Item(24,9).Value = 254
$ws.Cells.Item(24,10).Value = 266
$ws.Cells.Item(24,11).Value = -4.511278195488
$ws.Cells.Item(24,12).Value = -8.303249097472

# --- Row 25 ---
$ws.Cells.Item(25,3).Value = 7
$ws.Cells.Item(25,4).Value = 6
$ws.Cells.Item(25,4).NumberFormat = '#,##0'
$ws.Cells.Item(25,5).Value = 16.666666666666
$ws.Cells.Item(25,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(25,6).Value = 24
$ws.Cells.Item(25,7).Value = 16
$ws.Cells.Item(25,8).Value = 50
$ws.Cells.Item(25,9).Value = 68
$ws.Cells.Item(25,10).Value = 73
$ws.Cells.Item(25,11).Value = -6.849315068493
$ws.Cells.Item(25,12).Value = -30.612244897959

# --- Row 26 ---
$ws.Cells.Item(26,3).Value = 20
$ws.Cells.Item(26,4).Value = 16
$ws.Cells.Item(26,5).Value = 25
$ws.Cells.Item(26,6).Value = 64
$ws.Cells.Item(26,7).Value = 52
$ws.Cells.Item(26,8).Value = 23.076923076923
$ws.Cells.Item(26,9).Value = 192
$ws.Cells.Item(26,10).Value = 189
$ws.Cells.Item(26,11).Value = 1.587301587301
$ws.Cells.Item(26,12).Value = 37.142857142857
$ws.Cells.Item(26,13).Value = 68.421052631578

# --- Row 27 ---
$ws.Cells.Item(27,4).Value = "'0"
$ws.Cells.Item(14, 4).Copy() | Out-Null
$ws.Cells.Item(27,4).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,5).Value = "***.*"
$ws.Cells.Item(14, 5).Copy() | Out-Null
$ws.Cells.Item(27,5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27,6).Value = 2
$ws.Cells.Item(27,8).Value = -33.333333333333
$ws.Cells.Item(27,9).Value = 15
$ws.Cells.Item(27,11).Value = 66.666666666666
$ws.Cells.Item(27,12).Value = 275

# --- Row 28 ---
$ws.Cells.Item(28,3).Value = "'0"
$ws.Cells.Item(14, 4).Copy() | Out-Null
$ws.Cells.Item(28,3).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28,4).Value = 1
$ws.Cells.Item(28,4).NumberFormat = '#,##0'
$ws.Cells.Item(28,5).Value = -100
$ws.Cells.Item(28,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(28,7).Value = 2
$ws.Cells.Item(28,8).Value = 350
$ws.Cells.Item(28,9).Value = 21
$ws.Cells.Item(28,10).Value = 10
$ws.Cells.Item(28,11).Value = 110
$ws.Cells.Item(28,12).Value = 31.25

# --- Row 29 ---
$ws.Cells.Item(29,3).Value = 1
$ws.Cells.Item(29,3).NumberFormat = '#,##0'
$ws.Cells.Item(29,4).Value = 1
$ws.Cells.Item(29,4).NumberFormat = '#,##0'
$ws.Cells.Item(29,5).Value = 0
$ws.Cells.Item(29,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(29,6).Value = 1
$ws.Cells.Item(29,6).NumberFormat = '#,##0'
$ws.Cells.Item(29,7).Value = 1
$ws.Cells.Item(29,7).NumberFormat = '#,##0'
$ws.Cells.Item(29,8).Value = 0
$ws.Cells.Item(29,8).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(29,9).Value = 3
$ws.Cells.Item(29,10).Value = 1
$ws.Cells.Item(29,10).NumberFormat = '#,##0'
$ws.Cells.Item(29,11).Value = 200
$ws.Cells.Item(29,11).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(29,12).Value = -25
$ws.Cells.Item(29,13).Value = 0
$ws.Cells.Item(29,14).Value = -86.363636363636

# --- Row 30 ---
$ws.Cells.Item(30,3).Value = 1
$ws.Cells.Item(30,3).NumberFormat = '#,##0'
$ws.Cells.Item(30,4).Value = 1
$ws.Cells.Item(30,4).NumberFormat = '#,##0'
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(30,5).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(30,6).Value = 1
$ws.Cells.Item(30,6).NumberFormat = '#,##0'
$ws.Cells.Item(30,7).Value = 1
$ws.Cells.Item(30,7).NumberFormat = '#,##0'
$ws.Cells.Item(30,8).Value = 0
$ws.Cells.Item(30,8).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(30,9).Value = 3
$ws.Cells.Item(30,10).Value = 1
$ws.Cells.Item(30,10).NumberFormat = '#,##0'
$ws.Cells.Item(30,11).Value = 200
$ws.Cells.Item(30,11).NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Cells.Item(30,12).Value = 200
$ws.Cells.Item(30,13).Value = 0
$ws.Cells.Item(30,14).Value = -84.210526315789

